# Re-index the "_UP" (ascending) profile sheets so their bin index column (A)
# is offset to continue from where the corresponding "_PROFILE"/"_DOWN"
# sheets for the same station already start (matching the newly-binned data
# that was prepended). Every other "_UP"/"_DOWN"/"_PROFILE" sheet pair was
# already re-processed; station 74 ("_UP") needs a +11 shift (one extra bin
# vs. its siblings) while every other station's "_UP" sheet needs +10.

$wb = $excel.ActiveWorkbook

$shifts = @{
    "74_UP"  = 11;
    "83_UP"  = 10;
    "95_UP"  = 10;
    "106_UP" = 10;
    "113_UP" = 10;
    "118_UP" = 10;
    "127_UP" = 10;
    "133_UP" = 10;
}

foreach ($sheetName in $shifts.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $shift = $shifts[$sheetName]

    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $oldVal = $cell.Value2
        if ($oldVal -ne $null) {
            $cell.Value = $oldVal + $shift
        }
    }
}
